# Update two-digit-by-one-digit division problems to the newly generated set.
$d = $word.ActiveDocument

$replacements = @(
    @("92÷6=", "84÷8="),
    @("18÷2=", "89÷4="),
    @("46÷9=", "77÷7="),
    @("26÷9=", "78÷8="),
    @("47÷5=", "85÷8="),
    @("70÷7=", "95÷6="),
    @("97÷7=", "48÷3="),
    @("68÷8=", "81÷9="),
    @("73÷4=", "97÷2="),
    @("43÷5=", "89÷3="),
    @("79÷8=", "75÷2="),
    @("44÷3=", "11÷6="),
    @("53÷7=", "23÷6="),
    @("75÷5=", "27÷8="),
    @("40÷2=", "63÷9="),
    @("53÷4=", "91÷5="),
    @("30÷2=", "67÷6="),
    @("34÷2=", "56÷4="),
    @("39÷4=", "79÷4="),
    @("75÷3=", "71÷4="),
    @("19÷3=", "90÷9="),
    @("74÷3=", "95÷9="),
    @("91÷8=", "36÷7="),
    @("11÷3=", "13÷8="),
    @("50÷3=", "19÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
